$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "[SIAI] 왜 쉬운 과정은 없나요? 왜 장학금이 없나요?"
$ws.Range("E9").Value = "https://pdsi.pabii.com/siai-easy-coursework-scholarship/#utm_source=rss&utm_medium=rss&utm_campaign=siai-easy-coursework-scholarship"

$ws.Range("D12").Value = "SciPy Korea 2023 컨퍼런스 티켓 판매가 시작되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2023/04/26/scipy-korea-2023-%ec%bb%a8%ed%8d%bc%eb%9f%b0%ec%8a%a4-%ed%8b%b0%ec%bc%93-%ed%8c%90%eb%a7%a4%ea%b0%80-%ec%8b%9c%ec%9e%91%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"
